# Abänderungen und Zusätze zur Rückschau
# Target sheet: Tabelle1 (first worksheet) of TestprotokolleAlle.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
[void]$ws.Activate()

# Rows 8 and 10-19 in column C were attributed to "Tim Braumann"; they are
# now re-attributed to "Andre Kamp" (existing shared string). Row 9's
# (previously empty) C cell is filled in with the same tester name.
$testerRows = @(8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19)
foreach ($r in $testerRows) {
    $ws.Cells.Item($r, 3).Value = "Andre Kamp"
}

# Row 29's tester cell (C29) keeps its value ("Tim Braumann") but picks up
# the same font formatting already used by C27/C28 (explicit black font).
$ws.Range("C29").Font.Color = 0

# Move the visible viewport up a row and change the active selection.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C23").Select()
